$wb = $excel.ActiveWorkbook

# --- Tasks sheet: reformat the date columns (Due Date / Done Date / Add Date) ---
$tasks = $wb.Worksheets.Item("Tasks")
$tasks.Range("C2:E23").NumberFormat = "mm/dd/yy;@"

# --- Questions sheet: append two new question rows to the table ---
$questions = $wb.Worksheets.Item("Questions")
$lo = $questions.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$questions.Range("A5").Value = "What pace should I be going through these articles?"
$questions.Range("A6").Value = "How many articles should I cover in my survey presentation?"
$questions.Range("A7").Select() | Out-Null

# Re-select on the Tasks sheet last so it remains the active tab/sheet
$tasks.Range("A10").Select() | Out-Null
